$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7..55 down to 8..56
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new data record
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44602
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 100112032
$ws.Range("G7").Value = "Zapallo italiano"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 7500
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7750
$ws.Range("N7").Value = "$/caja 60 unidades"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 129
$ws.Range("Q7").Value = 60
$ws.Range("R7").Value = "Hortaliza"
